$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was HD HYUNDAI MIPO / 010620.KS) -> HDKSOE / 009540.KS
$ws.Range("B2").Value = "HDKSOE"
$ws.Range("C2").Value = "009540.KS"
$ws.Range("D2").Value = 426500
$ws.Range("E2").Value = 46.4
$ws.Range("F2").Value = 4.02
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 57.7
$ws.Range("N2").Value = 54.85170003294819

# Row 3 (was HDKSOE / 009540.KS) -> HD HYUNDAI MIPO / 010620.KS
$ws.Range("B3").Value = "HD HYUNDAI MIPO"
$ws.Range("C3").Value = "010620.KS"
$ws.Range("D3").Value = 223000
$ws.Range("E3").Value = 41.6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 56.5
$ws.Range("N3").Value = 54.85170003294819

# Row 4 (was SamsungHvyInd / 010140.KS) -> Hanwha Ocean / 042660.KS
$ws.Range("B4").Value = "Hanwha Ocean"
$ws.Range("C4").Value = "042660.KS"
$ws.Range("D4").Value = 107100
$ws.Range("E4").Value = 18
$ws.Range("F4").Value = -0.65
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 73
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 51.5
$ws.Range("N4").Value = 54.85170003294819

# Row 5 (was Hanwha Ocean / 042660.KS) -> SamsungHvyInd / 010140.KS
$ws.Range("B5").Value = "SamsungHvyInd"
$ws.Range("C5").Value = "010140.KS"
$ws.Range("D5").Value = 24900
$ws.Range("E5").Value = 39.5
$ws.Range("F5").Value = 1.22
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 63
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 93
$ws.Range("K5").Value = 47.7
$ws.Range("N5").Value = 54.85170003294819
